$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dados")

$ws.Range("A2").Value = 4.33631436872199
$ws.Range("B2").Value = 2.559828710052578

$ws.Range("A3").Value = 3.104068870015793
$ws.Range("B3").Value = 1.1354757379688

$ws.Range("A4").Value = 1.062341570308431
$ws.Range("B4").Value = 1.435495356548319

$ws.Range("A5").Value = 1.653148343425738
$ws.Range("B5").Value = 2.786676407812485

$ws.Range("A6").Value = 2.738549572119815
$ws.Range("B6").Value = 4.512419884876583
